# Auto-generated Excel COM-interop script
# Applies literal value updates to market-price derived columns (H-N)
# across multiple worksheets, per the scheduled-runner data refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 248.04762
$ws.Range("I6").Value = 172.61111
$ws.Range("J6").Value = 700.6667
$ws.Range("K6").Value = 517.8333299999999
$ws.Range("L6").Value = 2102.0001
$ws.Range("M6").Value = -405.8333299999999
$ws.Range("N6").Value = -2326.0001

$ws.Range("H74").Value = 3436.5186
$ws.Range("I74").Value = 3556.625
$ws.Range("J74").Value = 3261.818
$ws.Range("K74").Value = 3556.625
$ws.Range("L74").Value = 3261.818
$ws.Range("M74").Value = -2620.625
$ws.Range("N74").Value = -5133.818

$ws.Range("H76").Value = 4793976.5
$ws.Range("I76").Value = 7624177.5
$ws.Range("J76").Value = 4406.4614
$ws.Range("K76").Value = 7624177.5
$ws.Range("L76").Value = 4406.4614
$ws.Range("M76").Value = -7623862.5
$ws.Range("N76").Value = -5036.4614

$ws.Range("H77").Value = 3436.5186
$ws.Range("I77").Value = 3556.625
$ws.Range("J77").Value = 3261.818
$ws.Range("K77").Value = 17783.125
$ws.Range("L77").Value = 16309.09
$ws.Range("M77").Value = -13103.125
$ws.Range("N77").Value = -25669.09

$ws.Range("H79").Value = 4793976.5
$ws.Range("I79").Value = 7624177.5
$ws.Range("J79").Value = 4406.4614
$ws.Range("K79").Value = 7624177.5
$ws.Range("L79").Value = 4406.4614
$ws.Range("M79").Value = -7623085.5
$ws.Range("N79").Value = -6590.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1950
$ws.Range("I3").Value = 1950
$ws.Range("K3").Value = 1950
$ws.Range("M3").Value = -1835

$ws.Range("H74").Value = 1407.7858
$ws.Range("I74").Value = 1052.2703
$ws.Range("J74").Value = 2100.1052
$ws.Range("K74").Value = 1052.2703
$ws.Range("L74").Value = 2100.1052
$ws.Range("M74").Value = -178.2702999999999
$ws.Range("N74").Value = -3848.1052

$ws.Range("H77").Value = 1407.7858
$ws.Range("I77").Value = 1052.2703
$ws.Range("J77").Value = 2100.1052
$ws.Range("K77").Value = 5261.3515
$ws.Range("L77").Value = 10500.526
$ws.Range("M77").Value = -893.3514999999998
$ws.Range("N77").Value = -19236.526

$ws.Range("H122").Value = 1426657.5
$ws.Range("I122").Value = 1711435.6
$ws.Range("K122").Value = 5134306.800000001
$ws.Range("M122").Value = -5131856.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83334680
$ws.Range("I99").Value = 200000880
$ws.Range("J99").Value = 1684.1428
$ws.Range("K99").Value = 200000880
$ws.Range("L99").Value = 1684.1428
$ws.Range("M99").Value = -199999382
$ws.Range("N99").Value = -4680.1428

$ws.Range("H134").Value = 7249.952
$ws.Range("I134").Value = 9509.643
$ws.Range("J134").Value = 2730.5715
$ws.Range("K134").Value = 28528.929
$ws.Range("L134").Value = 8191.7145
$ws.Range("M134").Value = -25993.929
$ws.Range("N134").Value = -13261.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4526975
$ws.Range("I16").Value = 9617133
$ws.Range("J16").Value = 2390.3333
$ws.Range("K16").Value = 9617133
$ws.Range("L16").Value = 2390.3333
$ws.Range("M16").Value = -9616846
$ws.Range("N16").Value = -2964.3333

$ws.Range("H58").Value = 1346.3636
$ws.Range("I58").Value = 1145.5676
$ws.Range("J58").Value = 1759.1111
$ws.Range("K58").Value = 1145.5676
$ws.Range("L58").Value = 1759.1111
$ws.Range("M58").Value = -942.5676000000001
$ws.Range("N58").Value = -2165.1111

$ws.Range("H107").Value = 718.8421
$ws.Range("I107").Value = 359.8889
$ws.Range("J107").Value = 1041.9
$ws.Range("K107").Value = 359.8889
$ws.Range("L107").Value = 1041.9
$ws.Range("M107").Value = 1560.1111
$ws.Range("N107").Value = -4881.9

$ws.Range("H113").Value = 4526975
$ws.Range("I113").Value = 9617133
$ws.Range("J113").Value = 2390.3333
$ws.Range("K113").Value = 9617133
$ws.Range("L113").Value = 2390.3333
$ws.Range("M113").Value = -9614963
$ws.Range("N113").Value = -6730.3333

$ws.Range("H122").Value = 2528625.8
$ws.Range("I122").Value = 3970009.8
$ws.Range("J122").Value = 6203.5
$ws.Range("K122").Value = 11910029.4
$ws.Range("L122").Value = 18610.5
$ws.Range("M122").Value = -11907579.4
$ws.Range("N122").Value = -23510.5

$ws.Range("H134").Value = 3355.0938
$ws.Range("I134").Value = 3301.9614
$ws.Range("K134").Value = 9905.8842
$ws.Range("M134").Value = -7370.8842

$ws.Range("H136").Value = 1346.3636
$ws.Range("I136").Value = 1145.5676
$ws.Range("J136").Value = 1759.1111
$ws.Range("K136").Value = 3436.7028
$ws.Range("L136").Value = 5277.3333
$ws.Range("M136").Value = -886.7028
$ws.Range("N136").Value = -10377.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5865.6978
$ws.Range("I70").Value = 5922.5835
$ws.Range("J70").Value = 5573.143
$ws.Range("K70").Value = 5922.5835
$ws.Range("L70").Value = 5573.143
$ws.Range("M70").Value = -5652.5835
$ws.Range("N70").Value = -6113.143

$ws.Range("H73").Value = 5865.6978
$ws.Range("I73").Value = 5922.5835
$ws.Range("J73").Value = 5573.143
$ws.Range("K73").Value = 5922.5835
$ws.Range("L73").Value = 5573.143
$ws.Range("M73").Value = -4986.5835
$ws.Range("N73").Value = -7445.143

$ws.Range("H80").Value = 2474.36
$ws.Range("I80").Value = 2490.6667
$ws.Range("J80").Value = 2449.9
$ws.Range("K80").Value = 2490.6667
$ws.Range("L80").Value = 2449.9
$ws.Range("M80").Value = -1492.6667
$ws.Range("N80").Value = -4445.9

$ws.Range("H83").Value = 2474.36
$ws.Range("I83").Value = 2490.6667
$ws.Range("J83").Value = 2449.9
$ws.Range("K83").Value = 12453.3335
$ws.Range("L83").Value = 12249.5
$ws.Range("M83").Value = -7461.333500000001
$ws.Range("N83").Value = -22233.5

$ws.Range("H122").Value = 7939397.5
$ws.Range("I122").Value = 1325339
$ws.Range("J122").Value = 17860486
$ws.Range("K122").Value = 3976017
$ws.Range("L122").Value = 53581458
$ws.Range("M122").Value = -3973567
$ws.Range("N122").Value = -53586358

$ws.Range("H123").Value = 18748.637
$ws.Range("J123").Value = 19053.281
$ws.Range("L123").Value = 19053.281
$ws.Range("N123").Value = -23953.281

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 11344714
$ws.Range("I82").Value = 1250946.2
$ws.Range("K82").Value = 1250946.2
$ws.Range("M82").Value = -1250585.2

$ws.Range("H85").Value = 11344714
$ws.Range("I85").Value = 1250946.2
$ws.Range("K85").Value = 1250946.2
$ws.Range("M85").Value = -1249698.2

$ws.Range("H132").Value = 15158599
$ws.Range("I132").Value = 20841868
$ws.Range("J132").Value = 3216.5
$ws.Range("K132").Value = 62525604
$ws.Range("L132").Value = 9649.5
$ws.Range("M132").Value = -62523074
$ws.Range("N132").Value = -14709.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2039.8
$ws.Range("I122").Value = 1666.3334
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 4999.0002
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -2549.0002
$ws.Range("N122").Value = -12700

$ws.Range("H132").Value = 1801.5312
$ws.Range("I132").Value = 992.9524
$ws.Range("J132").Value = 3345.182
$ws.Range("K132").Value = 2978.8572
$ws.Range("L132").Value = 10035.546
$ws.Range("M132").Value = -448.8571999999999
$ws.Range("N132").Value = -15095.546

$ws.Range("H136").Value = 3231.7827
$ws.Range("I136").Value = 3614.742
$ws.Range("J136").Value = 2440.3333
$ws.Range("K136").Value = 10844.226
$ws.Range("L136").Value = 7320.999899999999
$ws.Range("M136").Value = -8294.226000000001
$ws.Range("N136").Value = -12420.9999
